$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(2)
try {
  $ser.Values(18) = 32
  Write-Host "indexer set ok"
} catch {
  Write-Host "indexer set failed: $_"
}
